$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 14:05"

# --- Simple numeric refreshes (country stays put) ---
# Iran (row 13)
$ws.Cells.Item(13,2).Value = 124603
$ws.Cells.Item(13,3).Value = 2111
$ws.Cells.Item(13,4).Value = 97173
$ws.Cells.Item(13,5).Value = 20311
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 62
$ws.Cells.Item(13,8).Value = 7119

# Paises Bajos (row 22)
$ws.Cells.Item(22,2).Value = 44249
$ws.Cells.Item(22,3).Value = 108
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = 21
$ws.Cells.Item(22,8).Value = 5715

# Catar (row 24)
$ws.Cells.Item(24,2).Value = 35606
$ws.Cells.Item(24,3).Value = 1637
$ws.Cells.Item(24,4).Value = 5634
$ws.Cells.Item(24,5).Value = 29957
$ws.Cells.Item(24,6).Value = 0
$ws.Cells.Item(24,7).Value = 0
$ws.Cells.Item(24,8).Value = 15

# Portugal (row 29)
$ws.Cells.Item(29,2).Value = 29432
$ws.Cells.Item(29,3).Value = 223
$ws.Cells.Item(29,4).Value = 6431
$ws.Cells.Item(29,5).Value = 21754
$ws.Cells.Item(29,6).Value = 0
$ws.Cells.Item(29,7).Value = 16
$ws.Cells.Item(29,8).Value = 1247

# --- Reorder Senegal / Sudan (rows 78-79) and refresh data ---
$ws.Cells.Item(78,1).Value = "Senegal"
$ws.Cells.Item(78,2).Value = 2617
$ws.Cells.Item(78,3).Value = 73
$ws.Cells.Item(78,4).Value = 1133
$ws.Cells.Item(78,5).Value = 1454
$ws.Cells.Item(78,6).Value = 0
$ws.Cells.Item(78,7).Value = 4
$ws.Cells.Item(78,8).Value = 30

$ws.Cells.Item(79,1).Value = "Sudan"
$ws.Cells.Item(79,2).Value = 2591
$ws.Cells.Item(79,3).Value = 0
$ws.Cells.Item(79,4).Value = 247
$ws.Cells.Item(79,5).Value = 2239
$ws.Cells.Item(79,6).Value = 0
$ws.Cells.Item(79,7).Value = 0
$ws.Cells.Item(79,8).Value = 105

# --- Reorder Libano / Albania (rows 107-108) and refresh data ---
$ws.Cells.Item(107,1).Value = "Libano"
$ws.Cells.Item(107,2).Value = 954
$ws.Cells.Item(107,3).Value = 23
$ws.Cells.Item(107,4).Value = 251
$ws.Cells.Item(107,5).Value = 677
$ws.Cells.Item(107,6).Value = 0
$ws.Cells.Item(107,7).Value = 0
$ws.Cells.Item(107,8).Value = 26

$ws.Cells.Item(108,1).Value = "Albania"
$ws.Cells.Item(108,2).Value = 949
$ws.Cells.Item(108,3).Value = 1
$ws.Cells.Item(108,4).Value = 742
$ws.Cells.Item(108,5).Value = 176
$ws.Cells.Item(108,6).Value = 0
$ws.Cells.Item(108,7).Value = 0
$ws.Cells.Item(108,8).Value = 31
